$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update CATEGORY (column D) for a few existing KEYWORDS rows:
# ESA, CNES, European Commission, UK Space Agency move from
# "Regulation & Policy" to the new "Institutional" category.
$ws.Range("D92").Value = "Institutional"
$ws.Range("D93").Value = "Institutional"
$ws.Range("D97").Value = "Institutional"
$ws.Range("D98").Value = "Institutional"

# --- Append new KEYWORDS / CATEGORY rows (135-143)
$newRows = @(
    @{ Row = 135; Keyword = "Space Armor"; Category = "Competitors" },
    @{ Row = 136; Keyword = "Debris shield"; Category = "Competitors" },
    @{ Row = 137; Keyword = "Atomic-6"; Category = "Competitors" },
    @{ Row = 138; Keyword = "Orbital debris"; Category = "Space Debris / Deorbit / ADR" },
    @{ Row = 139; Keyword = "France"; Category = "Institutional" },
    @{ Row = 140; Keyword = "IRIS2"; Category = "Satellite Operations" },
    @{ Row = 141; Keyword = "Hemeria"; Category = "Satellite Operations" },
    @{ Row = 142; Keyword = "French space command"; Category = "Institutional" },
    @{ Row = 143; Keyword = "Commandement de l'espace"; Category = "Institutional" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("C$r").Value = $item.Keyword
    $ws.Range("D$r").Value = $item.Category
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Refresh the view: show the newly added last row, selection on C143
$ws.Range("C143").Select()
